# feat : 마물, 보스 animator controller 데이터 추가
# Adds an "Anim" (Animator controller path) column to the Enemy and Boss
# data tables, mirroring the existing "Path" (model prefab) column.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Enemy sheet ("마물"/monster data) - append column G = Anim
# ---------------------------------------------------------------------
$enemy = $wb.Worksheets.Item("Enemy")

$enemy.Range("G1").Value = "Anim"
$enemy.Range("G2").Value = "string"
$enemy.Range("G3").Value = "Monster_Animatior"

$enemy.Range("G4").Value  = "Enemy/Enemy01Anim"
$enemy.Range("G6").Value  = "Enemy/Enemy03Anim"
$enemy.Range("G7").Value  = "Enemy/Enemy04Anim"
$enemy.Range("G8").Value  = "Enemy/Enemy05Anim"
$enemy.Range("G9").Value  = "Enemy/Enemy06Anim"
$enemy.Range("G10").Value = "Enemy/Enemy07Anim"
$enemy.Range("G11").Value = "Enemy/Enemy08Anim"
$enemy.Range("G12").Value = "Enemy/Enemy09Anim"
$enemy.Range("G13").Value = "Enemy/Enemy10Anim"
$enemy.Range("G5").Value  = "Enemy/Enemy02Anim"

# ---------------------------------------------------------------------
# Boss sheet ("보스" data) - append column G = Anim
# ---------------------------------------------------------------------
$boss = $wb.Worksheets.Item("Boss")

$boss.Range("G1").Value = "Anim"
$boss.Range("G2").Value = "string"
$boss.Range("G3").Value = "Monster_Animatior"

$boss.Range("G4").Value = "Boss/Boss01Anim"
$boss.Range("G6").Value = "Boss/Boss03Anim"
$boss.Range("G7").Value = "Boss/Boss04Anim"
$boss.Range("G8").Value = "Boss/Boss05Anim"
$boss.Range("G5").Value = "Boss/Boss02Anim"

# New print/page setup picked up on these two sheets in the same pass
# (both previously had no explicit pageSetup).
$chapter = $wb.Worksheets.Item("Chapter")
$chapter.PageSetup.PaperSize = 9
$chapter.PageSetup.Orientation = 1

$boss.PageSetup.PaperSize = 9
$boss.PageSetup.Orientation = 1

# ---------------------------------------------------------------------
# Selection / active-tab bookkeeping to match the author's final state:
# Unit -> Boss -> Enemy (Enemy ends up the active/selected sheet&cell).
# ---------------------------------------------------------------------
$unit = $wb.Worksheets.Item("Unit")
$unit.Activate()
$unit.Range("F15").Select()

$boss.Activate()
$boss.Range("G2").Select()

$enemy.Activate()
$enemy.Range("G2").Select()
